# Fix typo on slide 6b: "a significantly higher" -> "a slightly higher"
# (Correlations bullet about AWWA-policy compliant providers' percent water loss.)

$p = $ppt.ActivePresentation

$searchText = "a significantly "
$replacement = "a slightly "

$targetRange = $null

for ($si = 1; $si -le $p.Slides.Count; $si++) {
    $slide = $p.Slides.Item($si)
    for ($shi = 1; $shi -le $slide.Shapes.Count; $shi++) {
        $shape = $slide.Shapes.Item($shi)
        if (-not $shape.HasTextFrame) { continue }
        $tf = $shape.TextFrame
        if (-not $tf.HasText) { continue }

        $tr = $tf.TextRange
        $fullText = $tr.Text
        $idx0 = $fullText.IndexOf($searchText)

        if ($idx0 -ge 0) {
            # Grab just the sub-range that needs to change; this naturally
            # splits the parent run into the untouched-before / edited /
            # untouched-after runs, same as a user selecting the words and
            # retyping them in the PowerPoint UI.
            $startPos = $idx0 + 1   # PowerPoint ranges are 1-based
            $targetRange = $tr.Characters($startPos, $searchText.Length)
            break
        }
    }
    if ($targetRange -ne $null) { break }
}

if ($targetRange -eq $null) {
    throw "Could not locate the text '$searchText' to fix."
}

$targetRange.Text = $replacement
